$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 values
$row4 = @(42606.565162037034, -34, 61, 37, 14, 85, 11687, 7610, 432, 67, 41, 4, 23, "Named")
# New row 5 values
$row5 = @(42606.572326388887, -28, 68, 29, 11, 88, 14060, 9391, 498, 111, 48, 4, 30, "Named")

for ($c = 1; $c -le 14; $c++) {
    $ws.Cells.Item(4, $c).Value = $row4[$c - 1]
    $ws.Cells.Item(5, $c).Value = $row5[$c - 1]
}

# Apply the same date/time number format used in column A
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
